$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 913.85187
$ws.Range("J17").Value = 933.61536
$ws.Range("L17").Value = 2800.84608
$ws.Range("N17").Value = -3136.84608

$ws.Range("H19").Value = 879.4211
$ws.Range("I19").Value = 778.2
$ws.Range("J19").Value = 915.5714
$ws.Range("K19").Value = 778.2
$ws.Range("L19").Value = 915.5714
$ws.Range("M19").Value = -603.2
$ws.Range("N19").Value = -1265.5714

$ws.Range("H32").Value = 1700.7273
$ws.Range("J32").Value = 1986.7142
$ws.Range("L32").Value = 1986.7142
$ws.Range("N32").Value = -2638.7142

$ws.Range("H38").Value = 2383.7
$ws.Range("I38").Value = 2463.6
$ws.Range("J38").Value = 2303.8
$ws.Range("K38").Value = 7390.799999999999
$ws.Range("L38").Value = 6911.400000000001
$ws.Range("M38").Value = -7018.799999999999
$ws.Range("N38").Value = -7655.400000000001

$ws.Range("H96").Value = 558.5
$ws.Range("I96").Value = 517
$ws.Range("J96").Value = 600
$ws.Range("K96").Value = 1551
$ws.Range("L96").Value = 1800
$ws.Range("M96").Value = -178
$ws.Range("N96").Value = -4546

$ws.Range("H112").Value = 7353911.5
$ws.Range("J112").Value = 8621801
$ws.Range("L112").Value = 25865403
$ws.Range("N112").Value = -25867619

$ws.Range("H123").Value = 98926.336
$ws.Range("J123").Value = 98926.336
$ws.Range("L123").Value = 98926.336
$ws.Range("N123").Value = -108726.336

$ws.Range("H131").Value = 6245
$ws.Range("I131").Value = 2228.5
$ws.Range("J131").Value = 19633.334
$ws.Range("K131").Value = 6685.5
$ws.Range("L131").Value = 58900.00199999999
$ws.Range("M131").Value = -1645.5
$ws.Range("N131").Value = -68980.00199999999

$ws.Range("H137").Value = 1910.25
$ws.Range("I137").Value = 2160.0908
$ws.Range("J137").Value = 1360.6
$ws.Range("K137").Value = 6480.2724
$ws.Range("L137").Value = 4081.8
$ws.Range("M137").Value = -3930.2724
$ws.Range("N137").Value = -9181.799999999999

$ws.Range("H141").Value = 3759.7222
$ws.Range("I141").Value = 2519.4443
$ws.Range("K141").Value = 7558.3329
$ws.Range("M141").Value = -2378.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 991.85187
$ws.Range("I110").Value = 880.3182
$ws.Range("J110").Value = 1482.6
$ws.Range("K110").Value = 880.3182
$ws.Range("L110").Value = 1482.6
$ws.Range("M110").Value = 1164.6818
$ws.Range("N110").Value = -5572.6

$ws.Range("H123").Value = 33952.332
$ws.Range("J123").Value = 33952.332
$ws.Range("L123").Value = 33952.332
$ws.Range("N123").Value = -43752.332

$ws.Range("H130").Value = 35000
$ws.Range("J130").Value = 35000
$ws.Range("L130").Value = 35000
$ws.Range("N130").Value = -45040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7534.316
$ws.Range("I86").Value = 2277.3076
$ws.Range("J86").Value = 18924.5
$ws.Range("K86").Value = 2277.3076
$ws.Range("L86").Value = 18924.5
$ws.Range("M86").Value = -1154.3076
$ws.Range("N86").Value = -21170.5

$ws.Range("H89").Value = 7534.316
$ws.Range("I89").Value = 2277.3076
$ws.Range("J89").Value = 18924.5
$ws.Range("K89").Value = 11386.538
$ws.Range("L89").Value = 94622.5
$ws.Range("M89").Value = -5770.538
$ws.Range("N89").Value = -105854.5

$ws.Range("H94").Value = 947.4091
$ws.Range("I94").Value = 1010.15
$ws.Range("J94").Value = 320
$ws.Range("K94").Value = 1010.15
$ws.Range("L94").Value = 320
$ws.Range("M94").Value = -559.15
$ws.Range("N94").Value = -1222

$ws.Range("H105").Value = 2539.8157
$ws.Range("I105").Value = 2796.0454
$ws.Range("J105").Value = 2187.5
$ws.Range("K105").Value = 2796.0454
$ws.Range("L105").Value = 2187.5
$ws.Range("M105").Value = -1049.0454
$ws.Range("N105").Value = -5681.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 2000
$ws.Range("I25").Value = 2000
$ws.Range("K25").Value = 2000
$ws.Range("M25").Value = -1826

$ws.Range("H31").Value = 1553.0625
$ws.Range("I31").Value = 992.3333
$ws.Range("K31").Value = 992.3333
$ws.Range("M31").Value = -697.3333

$ws.Range("H34").Value = 1553.0625
$ws.Range("I34").Value = 992.3333
$ws.Range("K34").Value = 992.3333
$ws.Range("M34").Value = -790.3333

$ws.Range("H57").Value = 10074
$ws.Range("J57").Value = 10074
$ws.Range("L57").Value = 10074
$ws.Range("N57").Value = -11194

$ws.Range("H122").Value = 2027.5
$ws.Range("I122").Value = 1056
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 3168
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -718
$ws.Range("N122").Value = -13897

$ws.Range("H134").Value = 4171.174
$ws.Range("I134").Value = 2400.9285
$ws.Range("J134").Value = 6924.8887
$ws.Range("K134").Value = 7202.7855
$ws.Range("L134").Value = 20774.6661
$ws.Range("M134").Value = -4667.7855
$ws.Range("N134").Value = -25844.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 252.4
$ws.Range("I98").Value = 341
$ws.Range("J98").Value = 193.33333
$ws.Range("K98").Value = 1023
$ws.Range("L98").Value = 579.99999
$ws.Range("M98").Value = 475
$ws.Range("N98").Value = -3575.99999

$ws.Range("H131").Value = 962.1799999999999
$ws.Range("J131").Value = 1025.4131
$ws.Range("L131").Value = 3076.2393
$ws.Range("N131").Value = -13156.2393

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1237067.6
$ws.Range("J122").Value = 2875
$ws.Range("L122").Value = 8625
$ws.Range("N122").Value = -13525

$ws.Range("H126").Value = 2269.7942
$ws.Range("I126").Value = 1981.1666
$ws.Range("J126").Value = 2427.2273
$ws.Range("K126").Value = 5943.4998
$ws.Range("L126").Value = 7281.6819
$ws.Range("M126").Value = -3473.4998
$ws.Range("N126").Value = -12221.6819

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 15515
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 15515
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 15515
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -16105

$ws.Range("H27").Value = 15515
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 15515
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 15515
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -15729

$ws.Range("H36").Value = 22494.875
$ws.Range("J36").Value = 22494.875
$ws.Range("L36").Value = 22494.875
$ws.Range("N36").Value = -23618.875

$ws.Range("H128").Value = 53500
$ws.Range("J128").Value = 53500
$ws.Range("L128").Value = 53500
$ws.Range("N128").Value = -63460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 876.7857
$ws.Range("I107").Value = 858.875
$ws.Range("J107").Value = 900.6667
$ws.Range("K107").Value = 2576.625
$ws.Range("L107").Value = 2702.0001
$ws.Range("M107").Value = -656.625
$ws.Range("N107").Value = -6542.0001

$ws.Range("H122").Value = 43464.918
$ws.Range("I122").Value = 125894.75
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 377684.25
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -375234.25
$ws.Range("N122").Value = -11650

$ws.Range("H123").Value = 32140.2
$ws.Range("J123").Value = 32140.2
$ws.Range("L123").Value = 32140.2
$ws.Range("N123").Value = -41940.2

Write-Output "Edits applied successfully."
